$d = $word.ActiveDocument

# wdAlignParagraphLeft = 0, wdAlignParagraphRight = 2
$wdAlignParagraphRight = 2

$tbl = $d.Tables.Item(1)

# Header row: update text (case) and right-align the first column's cell.
$cell1 = $tbl.Cell(1, 1)
$cell1.Range.Text = "Numero competencia"
$cell1.Range.ParagraphFormat.Alignment = $wdAlignParagraphRight

$cell2 = $tbl.Cell(1, 2)
$cell2.Range.Text = "Text competencia"

# Data rows: right-align the first column (the competency number) only.
$rowCount = $tbl.Rows.Count
for ($i = 2; $i -le $rowCount; $i++) {
    $numCell = $tbl.Cell($i, 1)
    $numCell.Range.ParagraphFormat.Alignment = $wdAlignParagraphRight
}
